# Refresh the cryptocurrency price/volume table (GitHub Actions data pull).
# Rows keep their ranking position, but the coin occupying a given row,
# its price and its 1h volume change can all shift between pulls - so for
# every affected row we rewrite whichever of Coin/Link/Price/Volume cells
# changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the Price/Volume cells we are about to rewrite so
# that numeric-looking strings (e.g. "0.572", "7.45") are preserved verbatim
# as text instead of being auto-converted to floating point numbers.
$textCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "E25",
    "D26",
    "E26",
    "E27",
    "D28",
    "E28",
    "D29",
    "D30",
    "E30",
    "D31",
    "E31",
    "D32",
    "E32",
    "D33",
    "E33",
    "E34",
    "E35",
    "D36",
    "E36",
    "D37",
    "E37",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "E41",
    "E42",
    "D43",
    "E43",
    "E44",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values scraped from the latest cryptos feed.
$ws.Range("D2").Value = '42.415.29'
$ws.Range("E2").Value = '  -1.17%  '
$ws.Range("D3").Value = '2.520.54'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '316.41'
$ws.Range("E5").Value = '  +3.55%  '
$ws.Range("D6").Value = '94.13'
$ws.Range("E6").Value = '  -7.62%  '
$ws.Range("D7").Value = '0.572'
$ws.Range("E7").Value = '  -1.00%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").Value = '0.526'
$ws.Range("E9").Value = '  -3.89%  '
$ws.Range("D10").Value = '35.55'
$ws.Range("E10").Value = '  -6.43%  '
$ws.Range("D11").Value = '0.0804'
$ws.Range("E11").Value = '  -1.84%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.113'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '7.45'
$ws.Range("E13").Value = '  -4.01%  '
$ws.Range("D14").Value = '2.906.97'
$ws.Range("E14").Value = '  -0.82%  '
$ws.Range("D15").Value = '15.31'
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").Value = '2.471.58'
$ws.Range("E16").Value = '  -4.11%  '
$ws.Range("D17").Value = '0.839'
$ws.Range("E17").Value = '  -3.96%  '
$ws.Range("D18").Value = '42.458.11'
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("D19").Value = '12.98'
$ws.Range("E19").Value = '  -1.84%  '
$ws.Range("D20").Value = '6.54'
$ws.Range("E20").Value = '  +0.11%  '
$ws.Range("D21").Value = '0.0₃0955'
$ws.Range("E21").Value = '  -3.50%  '
$ws.Range("D22").Value = '69.89'
$ws.Range("E22").Value = '  -2.50%  '
$ws.Range("D23").Value = '250.09'
$ws.Range("E23").Value = '  -1.23%  '
$ws.Range("D24").Value = '2.95'
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("E25").Value = '  -2.91%  '
$ws.Range("D26").Value = '26.36'
$ws.Range("E26").Value = '  -3.42%  '
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("D28").Value = '2.38'
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("D29").Value = '10.12'
$ws.Range("D30").Value = '38.35'
$ws.Range("E30").Value = '  -1.89%  '
$ws.Range("D31").Value = '5.89'
$ws.Range("E31").Value = '  -5.24%  '
$ws.Range("D32").Value = '155.76'
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("D33").Value = '18.98'
$ws.Range("E33").Value = '  +3.04%  '
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("E35").Value = '  -0.21%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.0780'
$ws.Range("E36").Value = '  -2.39%  '
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '2.62'
$ws.Range("E37").Value = '  -1.26%  '
$ws.Range("E38").Value = '  -4.77%  '
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").Value = '23.88'
$ws.Range("E39").Value = '  -1.54%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.118'
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("E41").Value = '  +10.66%  '
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("D43").Value = '3.77'
$ws.Range("E43").Value = '  -3.42%  '
$ws.Range("E44").Value = '  -2.25%  '
$ws.Range("E45").Value = '  -6.27%  '
$ws.Range("D46").Value = '2.004.80'
$ws.Range("E46").Value = '  -2.83%  '
$ws.Range("D47").Value = '84.36'
$ws.Range("E47").Value = '  -2.23%  '
$ws.Range("D48").Value = '8.80'
$ws.Range("E48").Value = '  -2.52%  '
$ws.Range("D49").Value = '2.761.87'
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("D50").Value = '73.00'
$ws.Range("E50").Value = '  -0.79%  '
$ws.Range("D51").Value = '101.72'
$ws.Range("E51").Value = '  -1.50%  '
